$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Stage the new cell text ("1") in a scratch cell so it is stored as a
# *text* shared-string (Excel would otherwise infer "1" typed straight
# into B11 as a Number). Number-format it as Text first so the value
# survives round-tripping as a string.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "1"

# Copy only the value/type (not formatting) into B11 so its existing
# style (borders/fill/font) is preserved exactly as before.
$ws.Range("Z1").Copy()
$ws.Range("B11").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Clean up the scratch column so the sheet's used range/dimension is
# not left expanded by the helper cell.
$ws.Range("Z1").EntireColumn.Delete()
